# RMI file updates (again) 12/4
# - Recalibrate several "Share of New" values on the -psgr and -frgt sheets
# - Re-style the header row (bold + wrap on col A, right-align + wrap on cols B:H, taller row)
# - Re-point the active tab back to the "About" sheet / tidy up stale selections

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$psgr  = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$frgt  = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

function Update-HeaderRow($ws) {
    # Column A header: bold font + wrap text
    $colA = $ws.Range("A1")
    $colA.WrapText = $true
    $colA.Font.Bold = $true

    # Columns B:H header: right aligned + wrap text
    $colBH = $ws.Range("B1:H1")
    $colBH.WrapText = $true
    $colBH.HorizontalAlignment = -4152   # xlRight

    # Taller header row
    $ws.Rows.Item(1).RowHeight = 30
}

# ----- SoCDTtiNTY-psgr sheet -----
Update-HeaderRow $psgr

# Row 3 = HDVs: 0.044 -> 0.0435
$psgr.Range("B3:H3").Value2 = 0.0435
# Row 4 = aircraft: 0.046 -> 0.0416
$psgr.Range("B4:H4").Value2 = 0.0416
# Row 7 = motorbikes: 0.059 -> 0.0587
$psgr.Range("B7:H7").Value2 = 0.0587

# ----- SoCDTtiNTY-frgt sheet -----
Update-HeaderRow $frgt

# Row 3 = HDVs: 0.0353 -> 0.035
$frgt.Range("B3:H3").Value2 = 0.035

# ----- Window / selection bookkeeping -----
# Make sure every sheet's selection is reset to A1 (clears stale "A1:H7" selections)
$psgr.Activate()
$psgr.Range("A1").Select() | Out-Null
$frgt.Activate()
$frgt.Range("A1").Select() | Out-Null

# Re-activate the "About" sheet so it becomes the active tab again
$about.Activate()
$about.Range("A1").Select() | Out-Null
